# Weekly price update: a new price record for "Feria Lagunitas de Puerto
# Montt - Ajo" is inserted at row 104 (pushing the existing rows 104:149
# down to 105:150), matching the newest week's data point.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 104; existing rows 104-149 shift to 105-150.
$ws.Rows("104").Insert()

# Populate the new row 104 with the latest weekly record.
$ws.Range("A104").Value = 4
$ws.Range("B104").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C104").Value = "Los Lagos"
$ws.Range("D104").Value = 44466
$ws.Range("E104").Value = 10
$ws.Range("F104").Value = 100112003
$ws.Range("G104").Value = "Ajo"
$ws.Range("H104").Value = "Chino"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 180
$ws.Range("K104").Value = 17000
$ws.Range("L104").Value = 17000
$ws.Range("M104").Value = 17000
$ws.Range("N104").Value = "$/caja 10 kilos"
$ws.Range("O104").Value = "China"
$ws.Range("P104").Value = 1700
$ws.Range("Q104").Value = 10
$ws.Range("R104").Value = "Hortaliza"
